$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the data set. It belongs
# right before the existing row 45, so insert a blank row there first;
# this pushes the old rows 45..175 down to 46..176 (matching the diff's
# new dimension A1:R176).
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new observation. The
# descriptive columns (market/region/category/quality/unit/origin/etc.)
# repeat the same values used throughout this sheet.
$ws.Range("A45").Value2 = 11
$ws.Range("B45").Value = "Vega Monumental Concepción"
$ws.Range("C45").Value = "Bíobío"
$ws.Range("D45").Value2 = 44715
$ws.Range("D45").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E45").Value2 = 8
$ws.Range("F45").Value2 = 100112003
$ws.Range("G45").Value = "Ajo"
$ws.Range("H45").Value = "Chino"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value2 = 270
$ws.Range("K45").Value2 = 16000
$ws.Range("L45").Value2 = 17000
$ws.Range("M45").Value2 = 16556
$ws.Range("N45").Value = "$/caja 10 kilos"
$ws.Range("O45").Value = "China"
$ws.Range("P45").Value2 = 1656
$ws.Range("Q45").Value2 = 10
$ws.Range("R45").Value = "Hortaliza"
